$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 462.5
$ws.Range("I28").Value = 427.2857
$ws.Range("K28").Value = 427.2857
$ws.Range("M28").Value = 57.71429999999998

$ws.Range("H62").Value = 6258.5454
$ws.Range("I62").Value = 4818.8
$ws.Range("J62").Value = 7458.3335
$ws.Range("K62").Value = 4818.8
$ws.Range("L62").Value = 7458.3335
$ws.Range("M62").Value = -4194.8
$ws.Range("N62").Value = -8706.333500000001

$ws.Range("H65").Value = 6258.5454
$ws.Range("I65").Value = 4818.8
$ws.Range("J65").Value = 7458.3335
$ws.Range("K65").Value = 24094
$ws.Range("L65").Value = 37291.6675
$ws.Range("M65").Value = -20974
$ws.Range("N65").Value = -43531.6675

$ws.Range("H129").Value = 3141.2856
$ws.Range("I129").Value = 1962.3334
$ws.Range("J129").Value = 3462.818
$ws.Range("K129").Value = 5887.0002
$ws.Range("L129").Value = 10388.454
$ws.Range("M129").Value = -887.0002000000004
$ws.Range("N129").Value = -20388.454

$ws.Range("H135").Value = 709.2963
$ws.Range("I135").Value = 450.2381
$ws.Range("K135").Value = 4052.1429
$ws.Range("M135").Value = -1517.1429

$ws.Range("H137").Value = 2946.8708
$ws.Range("I137").Value = 1806
$ws.Range("J137").Value = 3770.8333
$ws.Range("K137").Value = 5418
$ws.Range("L137").Value = 11312.4999
$ws.Range("M137").Value = -2868
$ws.Range("N137").Value = -16412.4999

$ws.Range("H138").Value = 14068.7
$ws.Range("I138").Value = 7961.375
$ws.Range("J138").Value = 38498
$ws.Range("K138").Value = 23884.125
$ws.Range("L138").Value = 115494
$ws.Range("M138").Value = -18744.125
$ws.Range("N138").Value = -125774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13647.373
$ws.Range("I32").Value = 5847.8374
$ws.Range("K32").Value = 5847.8374
$ws.Range("M32").Value = -5560.8374

$ws.Range("H45").Value = 1399
$ws.Range("I45").Value = 1399
$ws.Range("K45").Value = 1399
$ws.Range("M45").Value = -1022

$ws.Range("H61").Value = 1423.0385
$ws.Range("I61").Value = 1425.96
$ws.Range("K61").Value = 1425.96
$ws.Range("M61").Value = -1213.96

$ws.Range("H136").Value = 1423.0385
$ws.Range("I136").Value = 1425.96
$ws.Range("K136").Value = 4277.88
$ws.Range("M136").Value = -1727.88

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 18000
$ws.Range("J76").Value = 18000
$ws.Range("L76").Value = 18000
$ws.Range("N76").Value = -18630

$ws.Range("H79").Value = 18000
$ws.Range("J79").Value = 18000
$ws.Range("L79").Value = 18000
$ws.Range("N79").Value = -20184

$ws.Range("H80").Value = 519.6
$ws.Range("I80").Value = 610
$ws.Range("J80").Value = 384
$ws.Range("K80").Value = 610
$ws.Range("L80").Value = 384
$ws.Range("M80").Value = 388
$ws.Range("N80").Value = -2380

$ws.Range("H83").Value = 519.6
$ws.Range("I83").Value = 610
$ws.Range("J83").Value = 384
$ws.Range("K83").Value = 3050
$ws.Range("L83").Value = 1920
$ws.Range("M83").Value = 1942
$ws.Range("N83").Value = -11904

$ws.Range("H105").Value = 4033.2058
$ws.Range("I105").Value = 3339.524
$ws.Range("J105").Value = 5153.769
$ws.Range("K105").Value = 3339.524
$ws.Range("L105").Value = 5153.769
$ws.Range("M105").Value = -1592.524
$ws.Range("N105").Value = -8647.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2221.5
$ws.Range("J16").Value = 2221.5
$ws.Range("L16").Value = 2221.5
$ws.Range("N16").Value = -2795.5

$ws.Range("H53").Value = 46008.332
$ws.Range("J53").Value = 46008.332
$ws.Range("L53").Value = 46008.332
$ws.Range("N53").Value = -47222.332

$ws.Range("H95").Value = 9812.5
$ws.Range("J95").Value = 9812.5
$ws.Range("L95").Value = 9812.5
$ws.Range("N95").Value = -15304.5

$ws.Range("H99").Value = 10285.634
$ws.Range("I99").Value = 7283.45
$ws.Range("K99").Value = 7283.45
$ws.Range("M99").Value = -5785.45

$ws.Range("H105").Value = 3613.6924
$ws.Range("I105").Value = 1921
$ws.Range("K105").Value = 1921
$ws.Range("M105").Value = -174

$ws.Range("H107").Value = 1051.9048
$ws.Range("I107").Value = 699.875
$ws.Range("K107").Value = 699.875
$ws.Range("M107").Value = 1220.125

$ws.Range("H113").Value = 2221.5
$ws.Range("J113").Value = 2221.5
$ws.Range("L113").Value = 2221.5
$ws.Range("N113").Value = -6561.5

$ws.Range("H126").Value = 10285.634
$ws.Range("I126").Value = 7283.45
$ws.Range("K126").Value = 21850.35
$ws.Range("M126").Value = -19380.35

$ws.Range("H132").Value = 3564.12
$ws.Range("I132").Value = 3171.95
$ws.Range("K132").Value = 9515.849999999999
$ws.Range("M132").Value = -6985.849999999999

$ws.Range("H134").Value = 4459.125
$ws.Range("I134").Value = 3645.8462
$ws.Range("J134").Value = 7983.3335
$ws.Range("K134").Value = 10937.5386
$ws.Range("L134").Value = 23950.0005
$ws.Range("M134").Value = -8402.5386
$ws.Range("N134").Value = -29020.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3380210
$ws.Range("I4").Value = 3380210
$ws.Range("K4").Value = 10140630
$ws.Range("M4").Value = -10140518

$ws.Range("H80").Value = 5927
$ws.Range("I80").Value = 5854
$ws.Range("K80").Value = 17562
$ws.Range("M80").Value = -16626

$ws.Range("H83").Value = 5927
$ws.Range("I83").Value = 5854
$ws.Range("K83").Value = 52686
$ws.Range("M83").Value = -48006

$ws.Range("H120").Value = 14609.048
$ws.Range("J120").Value = 16000
$ws.Range("L120").Value = 48000
$ws.Range("N120").Value = -57676

$ws.Range("H131").Value = 1035.8889
$ws.Range("I131").Value = 752.7143
$ws.Range("J131").Value = 2027
$ws.Range("K131").Value = 2258.1429
$ws.Range("L131").Value = 6081
$ws.Range("M131").Value = 2781.8571
$ws.Range("N131").Value = -16161

$ws.Range("H141").Value = 20000
$ws.Range("I141").Value = 20000
$ws.Range("K141").Value = 60000
$ws.Range("M141").Value = -54820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3962.3572
$ws.Range("I126").Value = 2946.2
$ws.Range("K126").Value = 8838.599999999999
$ws.Range("M126").Value = -6368.599999999999

$ws.Range("H132").Value = 2721.25
$ws.Range("I132").Value = 2136.76
$ws.Range("J132").Value = 7592
$ws.Range("K132").Value = 6410.280000000001
$ws.Range("L132").Value = 22776
$ws.Range("M132").Value = -3880.280000000001
$ws.Range("N132").Value = -27836

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2613.074
$ws.Range("I46").Value = 1658.6666
$ws.Range("K46").Value = 1658.6666
$ws.Range("M46").Value = -1470.6666

$ws.Range("H93").Value = 1285.5714
$ws.Range("I93").Value = 1106.5714
$ws.Range("J93").Value = 1464.5714
$ws.Range("K93").Value = 1106.5714
$ws.Range("L93").Value = 1464.5714
$ws.Range("M93").Value = 141.4286
$ws.Range("N93").Value = -3960.5714

$ws.Range("H132").Value = 4677.2285
$ws.Range("I132").Value = 3248.1428
$ws.Range("K132").Value = 9744.428400000001
$ws.Range("M132").Value = -7214.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = $null

$ws.Range("H96").Value = 950.875
$ws.Range("I96").Value = 1090.75
$ws.Range("J96").Value = 811
$ws.Range("K96").Value = 1090.75
$ws.Range("L96").Value = 811
$ws.Range("M96").Value = 282.25
$ws.Range("N96").Value = -3557

$ws.Range("H122").Value = 2047.375
$ws.Range("I122").Value = 2047.375
$ws.Range("K122").Value = 6142.125
$ws.Range("M122").Value = -3692.125

$ws.Range("H136").Value = 53157.65
$ws.Range("I136").Value = 2766.7334
$ws.Range("J136").Value = 204330.4
$ws.Range("K136").Value = 8300.200199999999
$ws.Range("L136").Value = 612991.2
$ws.Range("M136").Value = -5750.200199999999
$ws.Range("N136").Value = -618091.2
